$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 10.31211433333333
$ws.Range("H2").Value = 30.936343
$ws.Range("I2").Value = 0.633340936097251
$ws.Range("J2").Value = 0.633340936097251
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 6.346253666666667
$ws.Range("N2").Value = 19.038761
$ws.Range("O2").Value = 0.9446330608455225
$ws.Range("P2").Value = 0.9446330608455226
$ws.Range("Q2").Value = 65.44329339900256
$ws.Range("R2").Value = 588.9896405910231
$ws.Range("S2").Value = 0.5982747870243146
$ws.Range("T2").Value = 0.5982747870243147
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 10.31211433333333
$ws.Range("H3").Value = 30.936343
$ws.Range("I3").Value = 0.633340936097251
$ws.Range("J3").Value = 0.633340936097251
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.3719673333333333
$ws.Range("N3").Value = 1.115902
$ws.Range("O3").Value = 0.05536693915447755
$ws.Range("P3").Value = 0.05536693915447755
$ws.Range("Q3").Value = 3.835769669598444
$ws.Range("R3").Value = 34.521927026386
$ws.Range("S3").Value = 0.03506614907293635
$ws.Range("T3").Value = 0.03506614907293635
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.103438
$ws.Range("H4").Value = 12.310314
$ws.Range("I4").Value = 0.2520215719230645
$ws.Range("J4").Value = 0.2520215719230645
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 6.346253666666667
$ws.Range("N4").Value = 19.038761
$ws.Range("O4").Value = 0.9446330608455225
$ws.Range("P4").Value = 0.9446330608455226
$ws.Range("Q4").Value = 26.04145845343933
$ws.Range("R4").Value = 234.373126080954
$ws.Range("S4").Value = 0.2380679088847844
$ws.Range("T4").Value = 0.2380679088847844
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 4.103438
$ws.Range("H5").Value = 12.310314
$ws.Range("I5").Value = 0.2520215719230645
$ws.Range("J5").Value = 0.2520215719230645
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.3719673333333333
$ws.Range("N5").Value = 1.115902
$ws.Range("O5").Value = 0.05536693915447755
$ws.Range("P5").Value = 0.05536693915447755
$ws.Range("Q5").Value = 1.526344890358666
$ws.Range("R5").Value = 13.737104013228
$ws.Range("S5").Value = 0.0139536630382801
$ws.Range("T5").Value = 0.0139536630382801
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.866538
$ws.Range("H6").Value = 5.599614
$ws.Range("I6").Value = 0.1146374919796846
$ws.Range("J6").Value = 0.1146374919796846
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 6.346253666666667
$ws.Range("N6").Value = 19.038761
$ws.Range("O6").Value = 0.9446330608455225
$ws.Range("P6").Value = 0.9446330608455226
$ws.Range("Q6").Value = 11.84552362647267
$ws.Range("R6").Value = 106.609712638254
$ws.Range("S6").Value = 0.1082903649364235
$ws.Range("T6").Value = 0.1082903649364235
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.866538
$ws.Range("H7").Value = 5.599614
$ws.Range("I7").Value = 0.1146374919796846
$ws.Range("J7").Value = 0.1146374919796846
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.3719673333333333
$ws.Range("N7").Value = 1.115902
$ws.Range("O7").Value = 0.05536693915447755
$ws.Range("P7").Value = 0.05536693915447755
$ws.Range("Q7").Value = 0.6942911624253333
$ws.Range("R7").Value = 6.248620461828
$ws.Range("S7").Value = 0.006347127043261105
$ws.Range("T7").Value = 0.006347127043261104
